$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D1").Value = "bitcode"

$row = 2
for ($i = 2; $i -le 17; $i++) {
    $ws.Cells.Item($row, 4).Value = $i
    $row++
}
for ($i = 22; $i -le 37; $i++) {
    $ws.Cells.Item($row, 4).Value = $i
    $row++
}

$ws.Range("D1").Select()
